# "Format angepasst, Report gekürzt"
# The Suisa report was shortened: the deduction for the first entry (C2)
# dropped from 40% to 20%. Selection ends up resting on the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 20
$ws.Range("C2").Select()
